$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Constant "template" fragments (same text in every data row, 2-30)
# These used to be JS-object-literal style (unquoted keys, single-quoted img path)
# and are now proper JSON-ish (quoted keys, double-quoted img path).
$colA = '{"id": '
$colC = ', "clase": '
$colE = ', "producto": '
$colG = ', "tipo": '
$colI = ', "desc": '
$colK = ', "precio": '
$colM = ','
$colN = ' "img": '
$colP = '},'

for ($r = 2; $r -le 30; $r++) {
    $ws.Range("A$r").Value = $colA
    $ws.Range("C$r").Value = $colC
    $ws.Range("E$r").Value = $colE
    $ws.Range("G$r").Value = $colG
    $ws.Range("I$r").Value = $colI
    $ws.Range("K$r").Value = $colK
    $ws.Range("M$r").Value = $colM
    $ws.Range("N$r").Value = $colN
    $ws.Range("P$r").Value = $colP
}

# img column (O): switch from 'img/productos/NN.jpg' to "img/productos/NN.jpg".
# Prefixing the value with a single leading apostrophe makes Excel treat it as a
# quote-prefixed text entry (matching the original cell's quotePrefix style, s="2"),
# stripping that leading apostrophe from the stored text.
# Rows 23-26 keep a stray literal leading apostrophe in the final text (an authentic
# typo in the source edit) - achieved by doubling the leading apostrophe so the first
# is consumed as the prefix marker and the second remains as literal text.
$imgValues = @{
    2  = "'" + '"img/productos/01.jpg"'
    3  = "'" + '"img/productos/02.jpg"'
    4  = "'" + '"img/productos/03.jpg"'
    5  = "'" + '"img/productos/04.jpg"'
    6  = "'" + '"img/productos/05.jpg"'
    7  = "'" + '"img/productos/06.jpg"'
    8  = "'" + '"img/productos/07.jpg"'
    9  = "'" + '"img/productos/08.jpg"'
    10 = "'" + '"img/productos/09.jpg"'
    11 = "'" + '"img/productos/10.jpg"'
    12 = "'" + '"img/productos/11.jpg"'
    13 = "'" + '"img/productos/12.jpg"'
    14 = "'" + '"img/productos/13.jpg"'
    15 = "'" + '"img/productos/14.jpg"'
    16 = "'" + '"img/productos/15.jpg"'
    17 = "'" + '"img/productos/16.jpg"'
    18 = "'" + '"img/productos/17.jpg"'
    19 = "'" + '"img/productos/18.jpg"'
    20 = "'" + '"img/productos/19.jpg"'
    21 = "'" + '"img/productos/20.jpg"'
    22 = "'" + '"img/productos/21.jpg"'
    23 = "''" + '"img/productos/22.jpg"'
    24 = "''" + '"img/productos/23.jpg"'
    25 = "''" + '"img/productos/24.jpg"'
    26 = "''" + '"img/productos/25.jpg"'
    27 = "'" + '"img/productos/26.jpg"'
    28 = "'" + '"img/productos/27.jpg"'
    29 = "'" + '"img/productos/28.jpg"'
    30 = "'" + '"img/productos/29.jpg"'
}

for ($r = 2; $r -le 30; $r++) {
    $ws.Range("O$r").Value = $imgValues[$r]
}

# Restore the selected cell that was active when the workbook was last saved.
$null = $ws.Range("Q21").Select()
